$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 2613
$ws.Range("F4").Value = 466
$ws.Range("F5").Value = 294
$ws.Range("F7").Value = 473
$ws.Range("F8").Value = 1213
$ws.Range("F9").Value = 555
$ws.Range("F10").Value = 302
$ws.Range("F11").Value = 122
$ws.Range("F13").Value = 5648
$ws.Range("F15").Value = 1741
$ws.Range("F16").Value = 4090
$ws.Range("F17").Value = 422
$ws.Range("F20").Value = 4734
$ws.Range("F21").Value = 6150
$ws.Range("F23").Value = 1047
$ws.Range("F24").Value = 681
$ws.Range("F25").Value = 3733
$ws.Range("F26").Value = 491
$ws.Range("F28").Value = 189
$ws.Range("F29").Value = 125
$ws.Range("F30").Value = 980
$ws.Range("F31").Value = 1393
$ws.Range("F32").Value = 461
$ws.Range("F33").Value = 539
$ws.Range("F34").Value = 1583
$ws.Range("F36").Value = 1691
$ws.Range("F37").Value = 181
$ws.Range("F39").Value = 1113
$ws.Range("F41").Value = 1341
$ws.Range("F42").Value = 623
$ws.Range("F43").Value = 93
$ws.Range("F44").Value = 3361
$ws.Range("F45").Value = 126
$ws.Range("F46").Value = 276
$ws.Range("F47").Value = 407
$ws.Range("F49").Value = 3875

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 1189
$ws.Range("F6").Value = 40
$ws.Range("F9").Value = 1
$ws.Range("F24").Value = 70

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 3804

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 3804
$ws.Range("F4").Value = 2613
$ws.Range("F6").Value = 466
$ws.Range("F7").Value = 294
$ws.Range("F8").Value = 1189
$ws.Range("F11").Value = 473
$ws.Range("F12").Value = 1213
$ws.Range("F13").Value = 555
$ws.Range("F14").Value = 302
$ws.Range("F15").Value = 122
$ws.Range("F18").Value = 1741
$ws.Range("F19").Value = 4734
$ws.Range("F21").Value = 1047
$ws.Range("F22").Value = 681
$ws.Range("F23").Value = 3733
$ws.Range("F24").Value = 491
$ws.Range("F26").Value = 189
$ws.Range("F27").Value = 125
$ws.Range("F28").Value = 980
$ws.Range("F29").Value = 1393
$ws.Range("F30").Value = 461
$ws.Range("F31").Value = 539
$ws.Range("F33").Value = 1583
$ws.Range("F35").Value = 1691
$ws.Range("F37").Value = 1113
$ws.Range("F39").Value = 623
$ws.Range("F41").Value = 93
$ws.Range("F42").Value = 70
$ws.Range("F43").Value = 3361
$ws.Range("F45").Value = 126
$ws.Range("F46").Value = 276
$ws.Range("F47").Value = 407
$ws.Range("F49").Value = 3875
